$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New variable names (rows 2..17) replacing old ones, with dummy placeholder values
$names = @(
    "av_int_no_involvement_by_population",
    "av_int_event_photovoltaic_panels",
    "av_crop_ha",
    "av_crop_yield_t_ha",
    "av_crop_profit_EUR_t",
    "av_energy_ha",
    "av_energy_yield_kwh_ha",
    "av_energy_profit_EUR_kwh",
    "av_int_cost_search_panels",
    "av_int_cost_search_location",
    "av_int_cost_photovoltaic_panels",
    "av_int_cost_ground_preparation",
    "av_int_cost_installation",
    "av_int_cost_training",
    "av_int_cost_reparation",
    "av_int_benefit_shade"
)

$row = 2
foreach ($name in $names) {
    $ws.Cells.Item($row, 1).Value = $name
    $ws.Cells.Item($row, 2).Value = 0
    $ws.Cells.Item($row, 3).Value = "NA"
    $ws.Cells.Item($row, 4).Value = 1
    $ws.Cells.Item($row, 5).Value = "posnorm"
    $ws.Cells.Item($row, 6).Value = $name
    $ws.Cells.Item($row, 7).Value = $name
    $row++
}

$ws.Range("A28").Select()
